$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "last"
$ws.Range("B3").Value = " "
$ws.Range("C3").Value = " "
$ws.Range("D3").Value = "first"
$ws.Range("E3").Value = "'3333"
$ws.Range("F3").Value = "'333"
$ws.Range("G3").Value = "HCM"
$ws.Range("H3").Value = "address updated"

$ws.Range("A3:H3").Style = "Normal"
